$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 577, shifting all existing rows 577..658 down
# to 578..659 (their contents/formatting move with them, so the old row 658
# data ends up correctly duplicated into the new row 659).
$ws.Rows("577:577").Insert()

# Populate the newly inserted row 577 with its data (same constant columns
# as the rest of the "Repollo" / "Feria Lagunitas de Puerto Montt" block,
# plus the new record's own values).
$ws.Range("A577").Value2 = 4
$ws.Range("B577").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C577").Value2 = "Los Lagos"
$ws.Range("D577").Value2 = 44984
$ws.Range("E577").Value2 = 10
$ws.Range("F577").Value2 = 100112006
$ws.Range("G577").Value2 = "Repollo"
$ws.Range("H577").Value2 = "Crespo record"
$ws.Range("I577").Value2 = "Primera"
$ws.Range("J577").Value2 = 250
$ws.Range("K577").Value2 = 1700
$ws.Range("L577").Value2 = 1800
$ws.Range("M577").Value2 = 1740
$ws.Range("N577").Value2 = "$/unidad"
$ws.Range("O577").Value2 = "Región Metropolitana"
$ws.Range("P577").Value2 = 1740
$ws.Range("Q577").Value2 = 1
$ws.Range("R577").Value2 = "Hortaliza"
